# Extensions to the method definitions
# - Rename O1/P1 headers, insert 3 new columns (Q:variants, R:unused1, S:unused2)
#   between the old "special2" (now "filter") column and the "Comments" column,
#   moving the existing "Comments" column (and its data) from Q to T.
# - Populate two new cells in the (new) "variants" column.
# - Extend the AutoFilter / used range to account for the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 49

# 1) Move the existing "Comments" column (Q, col 17) three columns over to
#    column T (col 20), for every data + header row, before we overwrite Q.
for ($r = 2; $r -le $lastRow; $r++) {
    $src = $ws.Cells.Item($r, 17)
    $val = $src.Value2
    if ($null -ne $val) {
        $ws.Cells.Item($r, 20).Value = $val
        $src.Value = $null
    }
}

# 2) Rename / set the header row for the affected columns.
$ws.Cells.Item(1, 15).Value = "stars"      # O1: special1 -> stars
$ws.Cells.Item(1, 16).Value = "filter"     # P1: special2 -> filter
$ws.Cells.Item(1, 17).Value = "variants"   # Q1: new column
$ws.Cells.Item(1, 18).Value = "unused1"    # R1: new column
$ws.Cells.Item(1, 19).Value = "unused2"    # S1: new column
$ws.Cells.Item(1, 20).Value = "Comments"   # T1: moved from Q1
$ws.Cells.Item(1, 20).Font.Bold = $true
$ws.Cells.Item(1, 20).NumberFormat = "@"

# Two trailing header-styled (but empty) cells, U1:V1, matching the header
# row's formatting even though they hold no text (mirrors the original
# trailing blank header cells that used to sit at R1:S1).
$ws.Cells.Item(1, 21).Font.Bold = $true
$ws.Cells.Item(1, 21).NumberFormat = "@"
$ws.Cells.Item(1, 22).Font.Bold = $true
$ws.Cells.Item(1, 22).NumberFormat = "@"

# 2b) Row 19 has a wrap-text formatted block running across K:P; extend that
#     same formatting across the freshly inserted Q:S gap so the row's
#     formatting stays contiguous (matching the vacated Comments cell style).
$ws.Cells.Item(19, 17).WrapText = $true
$ws.Cells.Item(19, 18).WrapText = $true
$ws.Cells.Item(19, 19).WrapText = $true

# 3) New data values in the "variants" column.
$ws.Cells.Item(2, 17).Value = "asta_2019"
$ws.Cells.Item(47, 17).Value = "asta_future"

# 4) Extend dimension-dependent ranges (autofilter / filter database) to the
#    new right-hand edge of the table (column V, 22, with header through row 47).
$ws.Range("A1:V47").AutoFilter()
